$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 5272
$ws1.Range("F9").Value = 412
$ws1.Range("F13").Value = 129
$ws1.Range("F16").Value = 354
$ws1.Range("F22").Value = 6033
$ws1.Range("F26").Value = 6794
$ws1.Range("F32").Value = 4457
$ws1.Range("F34").Value = 134
$ws1.Range("F36").Value = 1127
$ws1.Range("F37").Value = 101
$ws1.Range("F40").Value = 915
$ws1.Range("F41").Value = 1112

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1148

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1148
$ws4.Range("F10").Value = 5272
$ws4.Range("F12").Value = 412
$ws4.Range("F16").Value = 129
$ws4.Range("F19").Value = 354
$ws4.Range("F26").Value = 6033
$ws4.Range("F30").Value = 6794
$ws4.Range("F36").Value = 4457
$ws4.Range("F39").Value = 134
$ws4.Range("F41").Value = 1127
$ws4.Range("F42").Value = 101
$ws4.Range("F45").Value = 915
$ws4.Range("F46").Value = 1112
